# Split the O'Rear et al. 2021 citation run so that "UC Davis " is
# inserted in front of "Suisun Marsh Fish Study" as its own run, i.e.:
#   "O'Rear, T., J. Durand, and P. Moyle. 2021. Suisun Marsh Fish Study. ..."
# becomes three runs:
#   "O'Rear, T., J. Durand, and P. Moyle. 2021. " | "UC Davis " | "Suisun Marsh Fish Study. ..."
$d = $word.ActiveDocument

# Locate the exact insertion boundary (just before "Suisun Marsh Fish Study")
# by searching for the unique trailing portion of the citation text. This
# collapses/moves the Range to the found match and gives us precise offsets,
# regardless of paragraph numbering.
$target = $d.Content
$found = $target.Find.Execute(
    "Suisun Marsh Fish Study. https://watershed.ucdavis.edu/project/suisun-marsh-fish-study.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the Suisun Marsh Fish Study citation text"
}

# Insert "UC Davis " right before the found text, with revision tracking on.
# Tracked insertions are recorded as their own run (distinct from the
# surrounding, unmodified text) without requiring any character-formatting
# changes. Accepting just that revision afterwards folds it back into a
# plain run while leaving the rest of the document (and its other runs'
# rsid/formatting attributes) untouched.
$revBefore = $d.Revisions.Count

$d.TrackRevisions = $true
$insPoint = $d.Range($target.Start, $target.Start)
$insPoint.InsertBefore("UC Davis ")
$d.TrackRevisions = $false

# Accept just the revision(s) created by the insertion above (there should
# be exactly one new one), leaving any unrelated pre-existing revisions
# untouched.
for ($i = $d.Revisions.Count; $i -gt $revBefore; $i--) {
    $d.Revisions.Item($i).Accept()
}
